# Refactor quick resolver value
# - Row 5 (export_entity) becomes save_nlp (drops its args)
# - Row 6 (save_nlp) becomes export_entity and takes over the person.xlsx args
# - Row 7 (save_storage) is removed entirely
# - Selection moves to A7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "save_nlp"
$ws.Range("D5").Value = ""
$ws.Range("C6").Value = "export_entity"
$ws.Range("D6").Value = "file_path=person.xlsx, label=PERSON"

$ws.Rows.Item(7).Delete() | Out-Null

$ws.Range("A7").Select() | Out-Null
